# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp text update
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 09:35"

# Singapur (row 30): update totals / new cases / recuperados
$ws.Range("B30").Value = 28343
$ws.Range("C30").Value = 305
$ws.Range("E30").Value = 18981

# Chequia (row 51)
$ws.Range("B51").Value = 8480
$ws.Range("C51").Value = 5
$ws.Range("D51").Value = 5468
$ws.Range("E51").Value = 2714

# Armenia (row 66)
$ws.Range("B66").Value = 4823
$ws.Range("C66").Value = 351
$ws.Range("D66").Value = 2019
$ws.Range("E66").Value = 2743
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 61

# El Salvador rises above Republica de Yibuti and Mayotte (rows 96-98 re-sorted
# by total cases, descending). Row 96 now holds El Salvador's updated figures,
# while Yibuti and Mayotte shift down one row each, keeping their own values.
$ws.Range("A96").Value = "El Salvador"
$ws.Range("B96").Value = 1413
$ws.Range("C96").Value = 75
$ws.Range("D96").Value = 474
$ws.Range("E96").Value = 909
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = 30

$ws.Range("A97").Value = "Republica de Yibuti"
$ws.Range("B97").Value = 1401
$ws.Range("C97").Value = 0
$ws.Range("D97").Value = 972
$ws.Range("E97").Value = 425
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 4

$ws.Range("A98").Value = "Mayotte"
$ws.Range("B98").Value = 1342
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 627
$ws.Range("E98").Value = 697
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 18

# Georgia (row 120)
$ws.Range("B120").Value = 701
$ws.Range("C120").Value = 6
$ws.Range("D120").Value = 432
$ws.Range("E120").Value = 257
